# "9th Stab - Cosmetic Changes"
# Insert two new columns (C and D) into the ratings table, pushing the
# existing date/price-target column from C to E, and seed the two new
# columns with the same "UN" ticker placeholder that lives in column B.
# The header row gets two brand-new week labels ("Jun_15" / "Jun_17")
# while the former B1/C1 headers ("Jun_13" / "Jun_10") slide right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank columns at C:D - this shifts the old column C (and its
# column-width formatting) to column E, leaving B untouched.
$ws.Columns("C:D").Insert()

# Match the new columns' width to the table's existing 8.0-wide column.
$ws.Columns("C:D").ColumnWidth = 7.166666666666667

# Row 1 (header): the old B1 label moves into the new D1 slot; B1 and C1
# get the two brand new labels.
$oldB1 = $ws.Range("B1").Text
$ws.Range("C1").Value = "Jun_15"
$ws.Range("B1").Value = "Jun_17"
$ws.Range("D1").Value = $oldB1

# Rows 2-27: fill the two new columns with the same ticker text as
# column B on that row.
$lastRow = 27
for ($r = 2; $r -le $lastRow; $r++) {
    $bText = $ws.Cells.Item($r, 2).Text
    $ws.Cells.Item($r, 3).Value = $bText
    $ws.Cells.Item($r, 4).Value = $bText
}
